$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 541, pushing existing row 541 (and below) down to 542
$ws.Rows.Item(541).Insert()

# Fill the new row 541 with values
$ws.Cells.Item(541, 1).Value = 5
$ws.Cells.Item(541, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(541, 3).Value = "Maule"
$ws.Cells.Item(541, 4).Value = 45142
$ws.Cells.Item(541, 4).NumberFormat = $ws.Cells.Item(542, 4).NumberFormat
$ws.Cells.Item(541, 5).Value = 7
$ws.Cells.Item(541, 6).Value = 100112032
$ws.Cells.Item(541, 7).Value = "Zapallo italiano"
$ws.Cells.Item(541, 8).Value = "Sin especificar"
$ws.Cells.Item(541, 9).Value = "Primera"
$ws.Cells.Item(541, 10).Value = 200
$ws.Cells.Item(541, 11).Value = 14000
$ws.Cells.Item(541, 12).Value = 14000
$ws.Cells.Item(541, 13).Value = 14000
$ws.Cells.Item(541, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(541, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(541, 16).Value = 280
$ws.Cells.Item(541, 17).Value = 50
$ws.Cells.Item(541, 18).Value = "Hortaliza"
